$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:B2").Copy()
$ws.Range("A3:B3").PasteSpecial(-4122)
$ws.Range("A3").Value = "watching aswell"
$ws.Range("B3").Value = ".*.txt"
$ws.Rows.Item(3).RowHeight = $ws.Rows.Item(2).RowHeight
